$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J12").Value = "18 hr"
$ws.Range("I12").Value = "0.601 / 0.7047 / 3.546"
$ws.Range("H12").Value = "0.6777 / 0.8254 / 1.485"
$ws.Range("G12").Value = "Done!"

$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("E7").Select()
